# "added slide titles for all slides"
#
# The deck originally had a single title slide. Six more slides are
# appended after it (Title and Content layout == slideLayout2.xml),
# each carrying just a title; the content placeholder is left empty,
# matching the target deck. The final slide's title text is long
# enough that PowerPoint shrinks it to fit ("shrink text on overflow"),
# so its text frame is switched to that auto-size mode.

$p = $ppt.ActivePresentation

$titles = @(
    "How research actually works",
    "A project-oriented view",
    "Project management",
    "Project organization",
    "Directory structure"
)

$idx = 2
foreach ($title in $titles) {
    $s = $p.Slides.Add($idx, 2)
    $s.Shapes.Item(1).TextFrame.TextRange.Text = $title
    $idx = $idx + 1
}

# Final slide: "Exercise: download and unzip git repo"
$s = $p.Slides.Add($idx, 2)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.AutoSize = 2 ; # ppAutoSizeTextToFitShape -> shrink text on overflow
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "Exercise: download and unzip git repo"
# "git" is flagged by the spell checker, which splits it into its own run
$gitRange = $titleRange.Characters(30, 3)
$gitRange.Text = "git"
